$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$sh = $s.Shapes.Item(2)
$tf = $sh.TextFrame

# Turn on "Shrink text on overflow" for the subtitle placeholder, which
# serializes as <a:bodyPr><a:normAutofit/></a:bodyPr>
$tf.AutoSize = 2

# Append a new paragraph "John Knox" after the existing subtitle text,
# keeping the original two runs ("A r" / "eal-time group tracking app")
# intact.
$tr = $tf.TextRange
$tr.InsertAfter("`rJohn Knox")
